$wb = $excel.ActiveWorkbook

# --- Update shared/header strings (rename "Rural Areas (...)" -> "...") ---
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("B1").Value = "National Average"
$wsMeans.Range("C1").Value = "State Average"

$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("B1").Value = "National Average SD"
$wsSD.Range("C1").Value = "State Average SD"

# --- Update "Means" sheet data values ---
$wsMeans.Range("B2").Value = 72
$wsMeans.Range("C2").Value = 62
$wsMeans.Range("D2").Value = 71
$wsMeans.Range("E2").Value = 51
$wsMeans.Range("F2").Value = 56
$wsMeans.Range("G2").Value = 64
$wsMeans.Range("B3").Value = 13
$wsMeans.Range("C3").Value = 32
$wsMeans.Range("D3").Value = 26
$wsMeans.Range("E3").Value = 47
$wsMeans.Range("F3").Value = 41
$wsMeans.Range("G3").Value = 33
$wsMeans.Range("B4").Value = 15
$wsMeans.Range("C4").Value = 5.8
$wsMeans.Range("D4").Value = 2.7
$wsMeans.Range("E4").Value = 1.3
$wsMeans.Range("F4").Value = 2.5
$wsMeans.Range("G4").Value = 2.5
$wsMeans.Range("B5").Value = 18
$wsMeans.Range("C5").Value = 5.1
$wsMeans.Range("D5").Value = 1.6
$wsMeans.Range("E5").Value = 1.8
$wsMeans.Range("F5").Value = 1.1
$wsMeans.Range("G5").Value = 1.7
$wsMeans.Range("B6").Value = 71
$wsMeans.Range("C6").Value = 54
$wsMeans.Range("D6").Value = 27
$wsMeans.Range("E6").Value = 27
$wsMeans.Range("F6").Value = 33
$wsMeans.Range("G6").Value = 38
$wsMeans.Range("B7").Value = 7.3
$wsMeans.Range("C7").Value = 10
$wsMeans.Range("D7").Value = 11
$wsMeans.Range("E7").Value = 18
$wsMeans.Range("F7").Value = 20
$wsMeans.Range("G7").Value = 18
$wsMeans.Range("B8").Value = 5.8
$wsMeans.Range("C8").Value = 8.3
$wsMeans.Range("D8").Value = 9.8
$wsMeans.Range("E8").Value = 8.3
$wsMeans.Range("F8").Value = 7.5
$wsMeans.Range("G8").Value = 7.7
$wsMeans.Range("B9").Value = 29
$wsMeans.Range("C9").Value = 41
$wsMeans.Range("D9").Value = 40
$wsMeans.Range("B10").Value = 0.37
$wsMeans.Range("C10").Value = 0.45
$wsMeans.Range("D10").Value = 0.5

# --- Update "Standard Deviations" sheet data values ---
$wsSD.Range("B2").Value = 27
$wsSD.Range("C2").Value = 33
$wsSD.Range("D2").Value = 19
$wsSD.Range("E2").Value = 30
$wsSD.Range("F2").Value = 32
$wsSD.Range("G2").Value = 30
$wsSD.Range("B3").Value = 23
$wsSD.Range("C3").Value = 33
$wsSD.Range("D3").Value = 18
$wsSD.Range("E3").Value = 30
$wsSD.Range("F3").Value = 32
$wsSD.Range("G3").Value = 31
$wsSD.Range("B4").Value = 16
$wsSD.Range("C4").Value = 7.2
$wsSD.Range("D4").Value = 1.1
$wsSD.Range("E4").Value = 1.3
$wsSD.Range("F4").Value = 3
$wsSD.Range("G4").Value = 2.7
$wsSD.Range("B5").Value = 22
$wsSD.Range("C5").Value = 7.5
$wsSD.Range("D5").Value = 2.8
$wsSD.Range("E5").Value = 2.7
$wsSD.Range("F5").Value = 2
$wsSD.Range("G5").Value = 2.3
$wsSD.Range("B6").Value = 37
$wsSD.Range("C6").Value = 26
$wsSD.Range("D6").Value = 4.3
$wsSD.Range("E6").Value = 8
$wsSD.Range("F6").Value = 15
$wsSD.Range("G6").Value = 15
$wsSD.Range("B7").Value = 8.7
$wsSD.Range("C7").Value = 11
$wsSD.Range("D7").Value = 0.65
$wsSD.Range("E7").Value = 19
$wsSD.Range("F7").Value = 19
$wsSD.Range("G7").Value = 17
$wsSD.Range("B8").Value = 7.8
$wsSD.Range("C8").Value = 10
$wsSD.Range("D8").Value = 3.9
$wsSD.Range("E8").Value = 7.9
$wsSD.Range("F8").Value = 10
$wsSD.Range("G8").Value = 8.8
$wsSD.Range("B9").Value = 10
$wsSD.Range("C9").Value = 31
$wsSD.Range("D9").Value = 0
$wsSD.Range("C10").Value = 0.1
$wsSD.Range("D10").Value = 0
